$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New weekly price observation arrives, pushing the existing history down
# one row (row 81 -> 82, ..., row 109 -> 110) and inserting the new record
# at row 81.
$ws.Rows.Item(81).Insert()

$ws.Cells.Item(81, 1).Value = 11
$ws.Cells.Item(81, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(81, 3).Value = "Bíobío"
$ws.Cells.Item(81, 4).Value = 45007
$ws.Cells.Item(81, 5).Value = 8
$ws.Cells.Item(81, 6).Value = 100112012
$ws.Cells.Item(81, 7).Value = "Espinaca"
$ws.Cells.Item(81, 8).Value = "Sin especificar"
$ws.Cells.Item(81, 9).Value = "Primera"
$ws.Cells.Item(81, 10).Value = 50
$ws.Cells.Item(81, 11).Value = 6000
$ws.Cells.Item(81, 12).Value = 6500
$ws.Cells.Item(81, 13).Value = 6200
$ws.Cells.Item(81, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(81, 15).Value = "Región Metropolitana"
$ws.Cells.Item(81, 16).Value = 620
$ws.Cells.Item(81, 17).Value = 10
$ws.Cells.Item(81, 18).Value = "Hortaliza"
